$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1915.75
$ws.Range("I86").Value = 897.375
$ws.Range("K86").Value = 897.375
$ws.Range("M86").Value = 225.625
$ws.Range("H89").Value = 1915.75
$ws.Range("I89").Value = 897.375
$ws.Range("K89").Value = 4486.875
$ws.Range("M89").Value = 1129.125
$ws.Range("H127").Value = 1712
$ws.Range("I127").Value = 694
$ws.Range("K127").Value = 2082
$ws.Range("M127").Value = 2878
$ws.Range("H132").Value = 1565
$ws.Range("I132").Value = 1389.2273
$ws.Range("K132").Value = 4167.6819
$ws.Range("M132").Value = -1637.6819

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2374
$ws.Range("I74").Value = 2358.2222
$ws.Range("K74").Value = 2358.2222
$ws.Range("M74").Value = -1484.2222
$ws.Range("H77").Value = 2374
$ws.Range("I77").Value = 2358.2222
$ws.Range("K77").Value = 11791.111
$ws.Range("M77").Value = -7423.111000000001
$ws.Range("H122").Value = 880.8
$ws.Range("I122").Value = 1134.6666
$ws.Range("K122").Value = 3403.9998
$ws.Range("M122").Value = -953.9998000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 64995
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 64995
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 64995
$ws.Range("N129").Value = -74995
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 5017500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 5017500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 5017500
$ws.Range("N131").Value = -5027580
$ws.Range("H132").Value = 19531.736
$ws.Range("I132").Value = 14100.5
$ws.Range("J132").Value = 27859.633
$ws.Range("K132").Value = 42301.5
$ws.Range("L132").Value = 83578.899
$ws.Range("M132").Value = -39771.5
$ws.Range("N132").Value = -88638.899
$ws.Range("H133").Value = 46910.125
$ws.Range("I133").Value = 54662.5
$ws.Range("J133").Value = 44326
$ws.Range("K133").Value = 54662.5
$ws.Range("L133").Value = 44326
$ws.Range("M133").Value = -52132.5
$ws.Range("N133").Value = -49386
$ws.Range("H134").Value = 5049.6313
$ws.Range("I134").Value = 2817.3438
$ws.Range("J134").Value = 16955.166
$ws.Range("K134").Value = 8452.0314
$ws.Range("L134").Value = 50865.49800000001
$ws.Range("M134").Value = -5917.0314
$ws.Range("N134").Value = -55935.49800000001
$ws.Range("H135").Value = 99999
$ws.Range("I135").Value = 99999
$ws.Range("J135").Value = 99999
$ws.Range("K135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("M135").Value = -94929
$ws.Range("N135").Value = -110139
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 39999
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 39999
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 39999
$ws.Range("N138").Value = -50279
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 396197.16
$ws.Range("I141").Value = 300114
$ws.Range("J141").Value = 456249.12
$ws.Range("K141").Value = 300114
$ws.Range("L141").Value = 456249.12
$ws.Range("M141").Value = -294934
$ws.Range("N141").Value = -466609.12

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120,I120,J120,K120,L120,M120,N120").ClearContents()
$ws.Range("H121,I121,J121,K121,L121,M121,N121").ClearContents()
$ws.Range("H122,I122,J122,K122,L122,N122").ClearContents()
$ws.Range("H123,I123,J123,K123,L123,M123").ClearContents()
$ws.Range("H124,I124,J124,K124,L124,M124").ClearContents()
$ws.Range("H125,I125,J125,K125,L125,M125").ClearContents()
$ws.Range("H126,I126,J126,K126,L126,M126,N126").ClearContents()
$ws.Range("H127,I127,J127,K127,L127,N127").ClearContents()
$ws.Range("H128,I128,J128,K128,L128,M128").ClearContents()
$ws.Range("H129,I129,J129,K129,L129,M129,N129").ClearContents()
$ws.Range("H130,I130,J130,K130,L130").ClearContents()
$ws.Range("H131,I131,J131,K131,L131,M131,N131").ClearContents()
$ws.Range("H132,I132,J132,K132,L132,M132,N132").ClearContents()
$ws.Range("H133,I133,J133,K133,L133,M133").ClearContents()
$ws.Range("H134,I134,J134,K134,L134,M134").ClearContents()
$ws.Range("H136,I136,J136,K136,L136,M136").ClearContents()
$ws.Range("H137,I137,J137,K137,L137,M137,N137").ClearContents()
$ws.Range("H138,I138,J138,K138,L138,M138").ClearContents()
$ws.Range("H139,I139,J139,K139,L139,M139,N139").ClearContents()
$ws.Range("H140,I140,J140,K140,L140,M140").ClearContents()
$ws.Range("H141,I141,J141,K141,L141,M141").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2537.6738
$ws.Range("I136").Value = 2227.7942
$ws.Range("K136").Value = 6683.382599999999
$ws.Range("M136").Value = -4133.382599999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6287.125
$ws.Range("I126").Value = 6495.2607
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 19485.7821
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -17015.7821
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 11401.75
$ws.Range("I132").Value = 5505.6665
$ws.Range("K132").Value = 16516.9995
$ws.Range("M132").Value = -13986.9995
$ws.Range("H136").Value = 2096.5293
$ws.Range("I136").Value = 1364.8334
$ws.Range("J136").Value = 3852.6
$ws.Range("K136").Value = 4094.5002
$ws.Range("L136").Value = 11557.8
$ws.Range("M136").Value = -1544.5002
$ws.Range("N136").Value = -16657.8
